$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 18:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1574507
$ws.Range("C4").Value = 3924
$ws.Range("D4").Value = 361531
$ws.Range("E4").Value = 1119228
$ws.Range("G4").Value = 215
$ws.Range("H4").Value = 93748

# Brasil (row 7)
$ws.Range("B7").Value = 275087
$ws.Range("C7").Value = 3202
$ws.Range("E7").Value = 150172
$ws.Range("G7").Value = 138
$ws.Range("H7").Value = 18121

# Alemania (row 11)
$ws.Range("B11").Value = 178170
$ws.Range("C11").Value = 343
$ws.Range("E11").Value = 13057
$ws.Range("G11").Value = 20
$ws.Range("H11").Value = 8213

# India (row 14)
$ws.Range("B14").Value = 111601
$ws.Range("C14").Value = 5126
$ws.Range("E14").Value = 65105
$ws.Range("G14").Value = 124
$ws.Range("H14").Value = 3426

# Canada (row 17)
$ws.Range("B17").Value = 79503
$ws.Range("C17").Value = 391
$ws.Range("D17").Value = 40342
$ws.Range("E17").Value = 33205
$ws.Range("G17").Value = 44
$ws.Range("H17").Value = 5956

# Singapur (row 30)
$ws.Range("D30").Value = 11207
$ws.Range("E30").Value = 18135

# Polonia (row 34)
$ws.Range("B34").Value = 19739
$ws.Range("C34").Value = 471
$ws.Range("E34").Value = 10594
$ws.Range("G34").Value = 14
$ws.Range("H34").Value = 962

# Egipto (row 44)
$ws.Range("B44").Value = 14229
$ws.Range("C44").Value = 745
$ws.Range("D44").Value = 3994
$ws.Range("E44").Value = 9555
$ws.Range("G44").Value = 21
$ws.Range("H44").Value = 680

# Republica Dominicana (row 45)
$ws.Range("B45").Value = 13477
$ws.Range("C45").Value = 254
$ws.Range("D45").Value = 7142
$ws.Range("E45").Value = 5889
$ws.Range("G45").Value = 5
$ws.Range("H45").Value = 446

# Argelia (row 56)
$ws.Range("B56").Value = 7542
$ws.Range("C56").Value = 165
$ws.Range("D56").Value = 3968
$ws.Range("E56").Value = 3006
$ws.Range("G56").Value = 7
$ws.Range("H56").Value = 568

# Irak (row 69) - stays in place, values updated
$ws.Range("B69").Value = 3724
$ws.Range("C69").Value = 113
$ws.Range("D69").Value = 2438
$ws.Range("E69").Value = 1152
$ws.Range("G69").Value = 3
$ws.Range("H69").Value = 134

# Row 70 becomes Azerbaiyan (moved up, new data)
$ws.Range("A70").Value = "Azerbaiyan"
$ws.Range("B70").Value = 3631
$ws.Range("C70").Value = 113
$ws.Range("D70").Value = 2253
$ws.Range("E70").Value = 1335
$ws.Range("G70").Value = 2
$ws.Range("H70").Value = 43

# Row 71 becomes Hungria (shifted down, data unchanged)
$ws.Range("A71").Value = "Hungria"
$ws.Range("B71").Value = 3598
$ws.Range("C71").Value = 42
$ws.Range("D71").Value = 1454
$ws.Range("E71").Value = 1674
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 470

# Row 72 becomes Camerun (shifted down, data unchanged)
$ws.Range("A72").Value = "Camerun"
$ws.Range("B72").Value = 3529
$ws.Range("C72").Value = 0
$ws.Range("D72").Value = 1567
$ws.Range("E72").Value = 1822
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 140

# Republica de Chipre (row 110)
$ws.Range("B110").Value = 922
$ws.Range("C110").Value = 4
$ws.Range("E110").Value = 389

# Reunion (row 131)
$ws.Range("B131").Value = 447
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 411
$ws.Range("E131").Value = 35
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 1

# Isla de Man (row 141)
$ws.Range("B141").Value = 336
$ws.Range("C141").Value = 1
$ws.Range("D141").Value = 302
$ws.Range("E141").Value = 10

# Row 196 becomes Santa Lucia (swap with Belice)
$ws.Range("A196").Value = "Santa Lucia"
$ws.Range("D196").Value = 18
$ws.Range("H196").Value = 0

# Row 197 becomes Belice (swap with Santa Lucia)
$ws.Range("A197").Value = "Belice"
$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2

# Row 209 becomes Montserrat (swap with Groenlandia)
$ws.Range("A209").Value = "Montserrat"
$ws.Range("D209").Value = 10
$ws.Range("H209").Value = 1

# Row 210 becomes Groenlandia (swap with Montserrat)
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0
